$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1029094938796418
$ws.Range("H2").Value = 55.17619670945673
$ws.Range("I2").Value = -34.14705762158022
$ws.Range("G3").Value = 0.05046814156701202
$ws.Range("H3").Value = -57.32545990796116
$ws.Range("G4").Value = -0.2936424821544729
$ws.Range("H4").Value = -6.780889640570628
$ws.Range("G5").Value = -0.3384204055881737
$ws.Range("H5").Value = 15.18601529871754
$ws.Range("G6").Value = 0.208212554431696
$ws.Range("H6").Value = 5.611675769885228
$ws.Range("G7").Value = 0.2709083859036889
$ws.Range("H7").Value = 30.63236953804124
$ws.Range("G8").Value = 0.1113414873493352
$ws.Range("H8").Value = 9.269714579722296
$ws.Range("G9").Value = 0.1282515833933437
$ws.Range("H9").Value = 1.400633536657717
$ws.Range("G10").Value = 0.04269980056232787
$ws.Range("H10").Value = -30.49914340071808
$ws.Range("G11").Value = 0.02188060448003902
$ws.Range("H11").Value = -56.17757206815033
$ws.Range("G12").Value = 0.08758475106929371
$ws.Range("H12").Value = -5.382532458870395
$ws.Range("G13").Value = 0.09943092346378567
$ws.Range("H13").Value = 30.47355990903425
$ws.Range("G14").Value = 0.1948511117687122
$ws.Range("H14").Value = -13.76805855561316
$ws.Range("G15").Value = 0.2359293311603966
$ws.Range("H15").Value = -4.234910366710644
$ws.Range("G16").Value = 0.1179941084073845
$ws.Range("H16").Value = 3.735773662455097
$ws.Range("G17").Value = 0.1218323568263161
$ws.Range("H17").Value = -18.45681427445225
$ws.Range("G18").Value = -0.01299707087588766
$ws.Range("H18").Value = -45.18704967858154
$ws.Range("G19").Value = 0.02929939039742637
$ws.Range("H19").Value = 20.96294842931101
$ws.Range("G20").Value = 0.1371809647756392
$ws.Range("H20").Value = 61.27333046853126
$ws.Range("G21").Value = 0.1076300077624839
$ws.Range("H21").Value = 64.43618929880697
$ws.Range("G22").Value = 0.1966239169325799
$ws.Range("H22").Value = 2.642615281383849
$ws.Range("G23").Value = 0.2156192224543509
$ws.Range("H23").Value = -0.04084257688065751
$ws.Range("G24").Value = -0.04050496729622739
$ws.Range("H24").Value = -964.3604463666671
$ws.Range("G25").Value = -0.02155564624548765
$ws.Range("H25").Value = 7.316728204337776
$ws.Range("G26").Value = 0.1758872547555272
$ws.Range("H26").Value = -14.14535410550951
$ws.Range("G27").Value = 0.1937367866898035
$ws.Range("H27").Value = 0.4418426891824712
$ws.Range("G28").Value = 0.04917498039396215
$ws.Range("H28").Value = -26.50981763570431
$ws.Range("G29").Value = 0.07370706085546821
$ws.Range("H29").Value = -21.80831574337999
